$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update text values
$ws.Range("B2").Value = "rohan"
$ws.Range("L2").Value = "Reading ,Drawing"
$ws.Range("B3").Value = "mini"
$ws.Range("H3").Value = "abcd"
$ws.Range("L3").Value = "Reading ,Writing"

# Update font color to explicit black (was theme color) for the Pincode/Phone
# cells that use the bordered font (I2, K2, I3, K3)
$ws.Range("I2").Font.Color = 0
$ws.Range("K2").Font.Color = 0
$ws.Range("I3").Font.Color = 0
$ws.Range("K3").Font.Color = 0

# Row heights
$ws.Rows.Item(1).RowHeight = 19.5
$ws.Rows.Item(2).RowHeight = 19.5
$ws.Rows.Item(3).RowHeight = 19.5
